# Costos Toma (Create, Update, Filter, Export)
#
# IndicadoresListado.xlsx:
#  - Remove the three workbook-level defined names (Indicadores_Mensual,
#    Indicadores_TablaM, Sucursales_Mensual).
#  - Replace the "{{Titulo}}" placeholder in each sheet's title cell (A1)
#    with a sheet-specific placeholder ({{TituloDiario}}, {{TituloSemanal}},
#    {{TituloMensual}}), keeping the existing rich-text formatting (bold/14,
#    plain/12, bold/12 runs).
#  - Clear the now-unused table header/body cells (A3, B3, A4) on every sheet.
#  - Move the active selection to A2 on every sheet (Mensual stays the tab
#    that is selected, matching the source file).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Drop all workbook-scoped defined names.
# ---------------------------------------------------------------------------
$nameCount = $wb.Names.Count
$i = $nameCount
while ($i -ge 1) {
    $wb.Names.Item($i).Delete()
    $i = $i - 1
}

# ---------------------------------------------------------------------------
# 2. Helper: rewrite the rich-text title cell, swapping the {{Titulo}} token
#    for a sheet-specific one while re-asserting the three original runs'
#    formatting (bold 14 / plain 12 / bold 12).
# ---------------------------------------------------------------------------
function Set-RichTitle {
    param($cell, $placeholder)

    $prefix = "Laboratorio Alfonso Ramos S.A. de C.V. MONTERREY`n{{Direccion}}`n{{Sucursal}}`n"
    $middle = "`n"
    $suffix = "Listado de " + $placeholder + "`nDEL {{Fecha}} AL {{Fecha}}"

    $cell.Value = $prefix + $middle + $suffix

    $len1 = $prefix.Length
    $len2 = $middle.Length
    $len3 = $suffix.Length

    $run1 = $cell.Characters(1, $len1)
    $run1.Font.Bold = $true
    $run1.Font.Size = 14

    $run2 = $cell.Characters($len1 + 1, $len2)
    $run2.Font.Bold = $false
    $run2.Font.Size = 12

    $run3 = $cell.Characters($len1 + $len2 + 1, $len3)
    $run3.Font.Bold = $true
    $run3.Font.Size = 12
}

# ---------------------------------------------------------------------------
# 3. Apply to each sheet: new title token, clear stale table cells, reset
#    selection to A2.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Diario")
Set-RichTitle $ws1.Range("A1") "{{TituloDiario}}"
$ws1.Range("A3").ClearContents()
$ws1.Range("B3").ClearContents()
$ws1.Range("A4").ClearContents()
$ws1.Range("A2").Select()

$ws2 = $wb.Worksheets.Item("Semanal")
Set-RichTitle $ws2.Range("A1") "{{TituloSemanal}}"
$ws2.Range("A3").ClearContents()
$ws2.Range("B3").ClearContents()
$ws2.Range("A4").ClearContents()
$ws2.Range("A2").Select()

$ws3 = $wb.Worksheets.Item("Mensual")
Set-RichTitle $ws3.Range("A1") "{{TituloMensual}}"
$ws3.Range("A3").ClearContents()
$ws3.Range("B3").ClearContents()
$ws3.Range("A4").ClearContents()
$ws3.Range("A2").Select()
